$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 2-3 values
$ws.Range("A2").Value = 3
$ws.Range("B2").Value = 6

$ws.Range("A3").Value = 0
$ws.Range("B3").Value = 5

# New rows 4-5 - copy style from A3 (bold/bordered/centered) so A4/A5 match formatting
$ws.Range("A3").Copy()
$ws.Range("A4:A5").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("A4").Value = 2
$ws.Range("B4").Value = 2

$ws.Range("A5").Value = 1
$ws.Range("B5").Value = 2
